# Fix the TODO placeholder in the speaker notes on slide 8 ("04 -EML
# standards.pptx"): replace the placeholder question with the final wording,
# and drop the trailing empty paragraph that followed it (the notes body
# ends up as a single paragraph of text).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$notesPage = $s.NotesPage

# Locate the notes body placeholder robustly (it currently holds the TODO
# text) rather than assuming a fixed shape index.
$notesShape = $null
for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
    $candidate = $notesPage.Shapes.Item($i)
    if ($candidate.HasTextFrame) {
        if ($candidate.TextFrame.TextRange.Text -like "*TODO*") {
            $notesShape = $candidate
        }
    }
}
if ($notesShape -eq $null) {
    $notesShape = $notesPage.Shapes.Item(1)
}

$apostrophe = [char]0x2019
$openQuote = [char]0x2018
$closeQuote = [char]0x2019

$newText = "So, we won" + $apostrophe + "t do this just yet, but it" + $apostrophe + "s worth highlighting that due to their shared understanding of EML we can take a dataset published on an IPT, and register it with a Pensoft journal as a " + $openQuote + "data paper" + $closeQuote + ", simply by exporting the EML and importing it to the appropriate Pensoft journal of our choice."

$notesShape.TextFrame.TextRange.Text = $newText
